# Apply "Atualizacoes 16 de janeiro de 2024" updates to the ValueSet workbook.

$wb = $excel.ActiveWorkbook

# --- Rename the second sheet tab ---
$includeSheet = $wb.Worksheets.Item("Include from ")
$includeSheet.Name = "Include from Estimated stroke"

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$meta.Range("B3").Value = "0.0.0"

# Title: shorter description
$meta.Range("B5").Value = "Estimated stroke severity using the NIHSS."

# Experimental value (was blank) -> false
$meta.Range("B7").Value = "false"

# Date updated
$meta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Description text tweak
$meta.Range("B12").Value = "ValueSet that intended to categorize the severity of the stroke based on the National Institutes of Health Stroke Scale (NIHSS)."

# --- Include sheet updates ---
$include = $wb.Worksheets.Item("Include from Estimated stroke")

# System URI value -> new CodeSystem URL
$include.Range("B9").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/EstStrokeSevNIHSSCatCS"
